# Updated cryptos list refresh: prices + hourly volume % changes, plus the
# Aave / WhiteBITCoin rows swapping places (rank index unchanged). Numeric-
# looking price strings are written with a leading apostrophe so Excel
# stores them as text (matching the source inlineStr cells) instead of
# silently parsing them into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.938.71'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = '2.509.14'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''532.87'
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = '''135.21'
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '2.953.56'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '58.838.09'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '''22.75'
$ws.Range("E15").Value = '  -2.27%  '
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '2.526.95'
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").Value = '''11.02'
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").Value = '''323.07'
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("E22").Value = '  +0.80%  '
$ws.Range("D23").Value = '''65.01'
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").Value = '''0.420'
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").Value = '0.0₃0763'
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("E29").Value = '  -3.66%  '
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").Value = '''169.19'
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("E33").Value = '  -4.86%  '
$ws.Range("E34").Value = '  -3.07%  '
$ws.Range("D35").Value = '''18.42'
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("E36").Value = '  -1.92%  '
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").Value = '''0.797'
$ws.Range("E39").Value = '  -4.53%  '
$ws.Range("D40").Value = '''281.02'
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").Value = '''0.605'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("E43").Value = '  -5.66%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").Value = '''10.93'
$ws.Range("E44").Value = '  +0.63%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''129.61'
$ws.Range("E45").Value = '  -1.03%  '
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("E47").Value = '  -2.54%  '
$ws.Range("E48").Value = '  -2.72%  '
$ws.Range("E49").Value = '  -0.90%  '
$ws.Range("D50").Value = '1.757.58'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("E51").Value = '  -0.51%  '
